# Regenerate save_data to use K (column G) instead of Strike# (old values).
# New K values were recomputed (regen std/mean, calc and write s_vals) and
# are written back into column G for rows 2-45, leaving the date/header row
# and the already-correct rows (40 and 42) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 2
    12 = 1
    13 = 1
    14 = 3
    15 = 0
    16 = 1
    17 = 2
    18 = 1
    19 = 3
    20 = 2
    21 = 3
    22 = 1
    23 = 1
    24 = 3
    25 = 0
    26 = 2
    27 = 0
    28 = 2
    29 = 1
    30 = 2
    31 = 3
    32 = 2
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 0
    39 = 2
    41 = 1
    43 = 2
    44 = 1
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
